# Updated cryptos list on Sun Feb 18 13:08:58 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while forcing it to stay as text (avoids Excel
# auto-converting numeric-looking strings such as "7.60" or "2.00" into
# numbers, which would silently drop meaningful trailing zeros), then
# restore the cell to its original "Normal" style so no visible formatting
# change is introduced.
function Set-TextValue {
    param($range, $text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "51.630.86"
$ws.Range("E2").Value = "  -0.06%  "
Set-TextValue $ws.Range("D3") "2.796.36"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("E4").Value = "  -0.05%  "
Set-TextValue $ws.Range("D5") "353.17"
$ws.Range("E5").Value = "  -1.69%  "
Set-TextValue $ws.Range("D6") "111.21"
$ws.Range("E6").Value = "  +2.38%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +7.12%  "
Set-TextValue $ws.Range("D10") "40.11"
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("E11").Value = "  -2.20%  "
Set-TextValue $ws.Range("D12") "0.0837"
$ws.Range("E12").Value = "  -1.03%  "
Set-TextValue $ws.Range("D13") "19.94"
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("E14").Value = "  +1.51%  "
Set-TextValue $ws.Range("D15") "3.235.03"
$ws.Range("E15").Value = "  +0.47%  "
Set-TextValue $ws.Range("D16") "2.822.15"
$ws.Range("E16").Value = "  +1.64%  "
Set-TextValue $ws.Range("D17") "0.945"
$ws.Range("E17").Value = "  +2.03%  "
Set-TextValue $ws.Range("D18") "51.562.67"
$ws.Range("E18").Value = "  -0.16%  "
Set-TextValue $ws.Range("D19") "7.60"
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("E20").Value = "  +3.28%  "
Set-TextValue $ws.Range("D21") "13.62"
$ws.Range("E21").Value = "  +3.22%  "
$ws.Range("E22").Value = "  +0.23%  "
Set-TextValue $ws.Range("D23") "70.25"
$ws.Range("E23").Value = "  +0.33%  "
Set-TextValue $ws.Range("D24") "266.88"
$ws.Range("E24").Value = "  -0.43%  "
Set-TextValue $ws.Range("D25") "2.75"
$ws.Range("E25").Value = "  -0.65%  "
$ws.Range("E26").Value = "  +0.08%  "
Set-TextValue $ws.Range("D27") "26.08"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("E28").Value = "  -2.52%  "
Set-TextValue $ws.Range("D29") "38.80"
$ws.Range("E29").Value = "  +9.83%  "
$ws.Range("E30").Value = "  +1.60%  "
$ws.Range("E31").Value = "  -3.92%  "
Set-TextValue $ws.Range("D32") "52.47"
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("E33").Value = "  -0.55%  "
$ws.Range("B34").Value = "VeChain"
$ws.Range("C34").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D34") "0.0450"
$ws.Range("E34").Value = "  +1.66%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D35") "0.0891"
$ws.Range("E35").Value = "  +6.09%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D36") "5.56"
$ws.Range("E36").Value = "  +7.77%  "
Set-TextValue $ws.Range("D37") "0.999"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("E38").Value = "  +0.04%  "
Set-TextValue $ws.Range("D39") "3.15"
$ws.Range("E39").Value = "  +0.50%  "
Set-TextValue $ws.Range("D40") "2.00"
$ws.Range("E40").Value = "  +2.44%  "
$ws.Range("E41").Value = "  +0.60%  "
Set-TextValue $ws.Range("D42") "2.50"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("E43").Value = "  +1.57%  "
Set-TextValue $ws.Range("D44") "121.15"
$ws.Range("E44").Value = "  +1.19%  "
Set-TextValue $ws.Range("D45") "21.73"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("E46").Value = "  +6.91%  "
$ws.Range("E47").Value = "  +4.34%  "
Set-TextValue $ws.Range("D48") "2.102.08"
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("E49").Value = "  +2.12%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue $ws.Range("D50") "62.25"
$ws.Range("E50").Value = "  +7.90%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws.Range("D51") "5.45"
$ws.Range("E51").Value = "  -1.69%  "
